$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.412.80'
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.160.13'
$ws.Range("E3").Value = '  +3.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.20'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("E6").Value = '  +1.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.26'
$ws.Range("E7").Value = '  +4.31%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.394'
$ws.Range("E9").Value = '  +2.19%  '

$ws.Range("E10").Value = '  +2.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.20'
$ws.Range("E12").Value = '  +8.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.480.94'
$ws.Range("E13").Value = '  +3.34%  '

$ws.Range("E14").Value = '  +2.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.818'
$ws.Range("E15").Value = '  +2.52%  '

$ws.Range("E16").Value = '  +1.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.144.55'
$ws.Range("E17").Value = '  +2.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.531.07'
$ws.Range("E18").Value = '  +2.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.45'
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.16'
$ws.Range("E20").Value = '  +1.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0855'
$ws.Range("E21").Value = '  +1.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.99'
$ws.Range("E22").Value = '  +0.87%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  +1.01%  '

$ws.Range("E25").Value = '  +0.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.79'
$ws.Range("E26").Value = '  +3.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.52'
$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("E28").Value = '  -0.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.66'
$ws.Range("E29").Value = '  +2.53%  '

$ws.Range("E30").Value = '  -2.42%  '

$ws.Range("E31").Value = '  +8.02%  '

$ws.Range("E32").Value = '  +1.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.65'
$ws.Range("E33").Value = '  +3.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.84'
$ws.Range("E34").Value = '  +2.93%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.19'
$ws.Range("E35").Value = '  +11.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0623'
$ws.Range("E36").Value = '  +1.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.42'
$ws.Range("E37").Value = '  +1.11%  '

$ws.Range("E38").Value = '  +0.54%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.29'
$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("E41").Value = '  +3.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.31'
$ws.Range("E42").Value = '  +2.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.538.22'
$ws.Range("E43").Value = '  -0.24%  '

$ws.Range("E44").Value = '  +6.40%  '

$ws.Range("E45").Value = '  +7.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0926'
$ws.Range("E46").Value = '  +0.49%  '

$ws.Range("E47").Value = '  -0.39%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.80'
$ws.Range("E48").Value = '  +1.96%  '

$ws.Range("E49").Value = '  +1.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.364.89'
$ws.Range("E50").Value = '  +3.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.97'
$ws.Range("E51").Value = '  -0.13%  '
